$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.129.04"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "1.867.86"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'307.05"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.5118"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D9").Value = "'0.07135"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("D10").Value = "'0.8879"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.914.72"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07528"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("D14").Value = "'5.317"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").Value = "'89.20"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'0.000008479"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").Value = "'14.12"
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("D19").Value = "'1.0000"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "27.177.49"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").Value = "2.071.88"
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("E23").Value = "  -2.74%  "
$ws.Range("D24").Value = "'6.469"
$ws.Range("E24").Value = "  -2.24%  "
$ws.Range("D25").Value = "'150.02"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").Value = "'1.843"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -2.77%  "
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("D29").Value = "'112.74"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "'4.748"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("D31").Value = "'4.672"
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("D32").Value = "'0.09026"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'0.05128"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("D34").Value = "'3.096"
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").Value = "'1.161"
$ws.Range("E35").Value = "  -6.51%  "
$ws.Range("D36").Value = "'0.7336"
$ws.Range("E36").Value = "  -6.28%  "
$ws.Range("D37").Value = "'0.02044"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "'2.502"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("D39").Value = "'3.052"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "'1.076"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").Value = "'6.603"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("D43").Value = "'116.76"
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").Value = "'0.1473"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").Value = "'0.9997"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'0.4624"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("E48").Value = "  -5.18%  "
$ws.Range("E49").Value = "  -4.29%  "
$ws.Range("D50").Value = "'64.51"
$ws.Range("E50").Value = "  -4.50%  "
$ws.Range("D51").Value = "'36.51"
$ws.Range("E51").Value = "  -1.69%  "
